$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- alpha_distance_range (row 2): Min/Max updated ---
$ws.Range("B2").Value = 5.6
$ws.Range("C2").Value = 10.4

# --- beta_distance_range (row 3): Min/Max updated ---
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 9.3

# --- ratio_threshold_range (row 4): left unchanged ---

# --- pie_threshold_range (currently row 6): Min/Max updated ---
# Update these values before deleting the theta_threshold_range row above it,
# so they land on the correct row once everything shifts up.
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 15

# Remove the theta_threshold_range row (row 5) entirely; pie_threshold_range
# shifts up to become the new row 5.
$ws.Rows.Item(5).Delete()

# Match the saved cursor position / print setup seen in the final workbook.
$null = $ws.Range("C4").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
